$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two additional data rows into the table (table grows from 16 to
#    18 data rows). Insert them just above the current last data row (31),
#    which pushes it down to become the new last row (33) and pushes the
#    footer block (36/37) down to (38/39) automatically.
#    Use PasteSpecial(xlPasteFormats) from a normal "middle" row (30) so the
#    new rows inherit the same cell styles (borders, fonts, number formats)
#    instead of a blank/default style.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

$ws.Rows("31:31").Insert()
$ws.Range("B30:J30").Copy()
$ws.Range("B31:J31").PasteSpecial($xlPasteFormats)

$ws.Rows("31:31").Insert()
$ws.Range("B30:J30").Copy()
$ws.Range("B31:J31").PasteSpecial($xlPasteFormats)

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Rewrite the worker/period detail table (rows 16-33) with the new data.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047442859"
$ws.Range("D16").Value = "ALVARO DE JESUS MALLARINO SANCHEZ"
$ws.Range("E16").Value = "2209"
$ws.Range("F16").Value = 57746
$ws.Range("G16").Value = 1883000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73200363"
$ws.Range("D17").Value = "WILLIAMS JOSE GUEVARA GOMEZ"
$ws.Range("E17").Value = "2501"
$ws.Range("F17").Value = 65520
$ws.Range("G17").Value = 1638000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1048434925"
$ws.Range("D18").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 52000
$ws.Range("G18").Value = 1300000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1048434925"
$ws.Range("D19").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 52000
$ws.Range("G19").Value = 1300000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1048434925"
$ws.Range("D20").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E20").Value = "2505"
$ws.Range("F20").Value = 52000
$ws.Range("G20").Value = 1300000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1048434925"
$ws.Range("D21").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E21").Value = "2504"
$ws.Range("F21").Value = 52000
$ws.Range("G21").Value = 1300000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1048434925"
$ws.Range("D22").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E22").Value = "2503"
$ws.Range("F22").Value = 52000
$ws.Range("G22").Value = 1300000

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1048434925"
$ws.Range("D23").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E23").Value = "2502"
$ws.Range("F23").Value = 52000
$ws.Range("G23").Value = 1300000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1048434925"
$ws.Range("D24").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E24").Value = "2501"
$ws.Range("F24").Value = 52000
$ws.Range("G24").Value = 1300000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1048434925"
$ws.Range("D25").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E25").Value = "2412"
$ws.Range("F25").Value = 52000
$ws.Range("G25").Value = 1300000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1048434925"
$ws.Range("D26").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E26").Value = "2411"
$ws.Range("F26").Value = 52000
$ws.Range("G26").Value = 1300000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1048434925"
$ws.Range("D27").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E27").Value = "2410"
$ws.Range("F27").Value = 52000
$ws.Range("G27").Value = 1300000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "1048434925"
$ws.Range("D28").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E28").Value = "2409"
$ws.Range("F28").Value = 52000
$ws.Range("G28").Value = 1300000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "1048434925"
$ws.Range("D29").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E29").Value = "2408"
$ws.Range("F29").Value = 52000
$ws.Range("G29").Value = 1300000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "1048434925"
$ws.Range("D30").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E30").Value = "2407"
$ws.Range("F30").Value = 52000
$ws.Range("G30").Value = 1300000

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "1048434925"
$ws.Range("D31").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E31").Value = "2406"
$ws.Range("F31").Value = 52000
$ws.Range("G31").Value = 1300000

$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "1048434925"
$ws.Range("D32").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E32").Value = "2405"
$ws.Range("F32").Value = 52000
$ws.Range("G32").Value = 1300000

$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "1048434925"
$ws.Range("D33").Value = "AURY LUZ GAVIRIA PUERTA"
$ws.Range("E33").Value = "2404"
$ws.Range("F33").Value = 27733
$ws.Range("G33").Value = 1300000

# ---------------------------------------------------------------------------
# 3. Update the summary figures above the table.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 930999
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 17
